# OLX Monitor update — 2026-02-17 12:01
# The monitoring run re-checked the same listings it saw at 11:35/11:36 and
# appended another identical batch of 8 rows (poqui x4, pokojewlublinie x2,
# dawnypatron x2) to the PODSUMOWANIE log, stamped with the new check time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

# Duplicate the most recent 8-row batch (rows 7-14) into rows 23-30,
# bringing along its values, number formats, and cell styles.
$srcRange = $ws.Range("A7:H14")
$dstRange = $ws.Range("A23:H30")
$srcRange.Copy($dstRange)

# Stamp the newly appended rows with this run's "last checked" timestamp.
$ws.Range("A23:A30").Value = "2026-02-17 12:01:39"
